$wb = $excel.ActiveWorkbook

$wsAssoc = $wb.Worksheets.Item("emshaps90")
$wsJob   = $wb.Worksheets.Item("NewJobInfo")

# --- emshaps90: widen column C (closest reachable snap to 18.77734375 chars) ---
$wsAssoc.Columns.Item(3).ColumnWidth = 17.92

# --- emshaps90: Associate Task rows (taskType / taskDate / taskNote get sample data) ---
# Written in the same order the original author must have used, so the new
# shared-string entries land at indices 48-53 in this exact sequence.

# Row 10 - taskDate sample values (entered as text, like the existing date-like
# entries in row 3, hence the leading apostrophe to force text storage)
$wsAssoc.Range("B10").Value = "'02/08/2017"
$wsAssoc.Range("C10").Value = "'02/14/2017"
$wsAssoc.Range("D10").Value = "STOP"
$wsAssoc.Range("A10").NumberFormat = "General"
$wsAssoc.Range("G10").NumberFormat = "General"

# Row 11 - taskNote sample values
$wsAssoc.Range("B11").Value = "Everything you want"
$wsAssoc.Range("C11").Value = "Everything you need"
$wsAssoc.Range("D11").Value = "STOP"
$wsAssoc.Range("A11").NumberFormat = "General"
$wsAssoc.Range("G11").NumberFormat = "General"

# Row 9 - taskType sample values
$wsAssoc.Range("B9").Value = "Certification"
$wsAssoc.Range("C9").Value = "Panel"
$wsAssoc.Range("D9").Value = "STOP"
$wsAssoc.Range("A9").NumberFormat = "General"
$wsAssoc.Range("G9").NumberFormat = "General"

# --- emshaps90: add a "Job Event" section (rows 13-18), mirroring the block
# already present on the NewJobInfo sheet so Associate Tasks can reference it ---
$wsAssoc.Range("A13").Value = "Job Event"
$wsAssoc.Range("A13").Font.Bold = $true

$wsAssoc.Range("A14").Value = "associate"
$wsAssoc.Range("A14").NumberFormat = "General"

$wsAssoc.Range("A15").Value = "assignment"
$wsAssoc.Range("A15").NumberFormat = "General"

$wsAssoc.Range("A16").Value = "eventType"
$wsAssoc.Range("A16").NumberFormat = "General"

$wsAssoc.Range("A17").Value = "eventDate"
$wsAssoc.Range("A17").NumberFormat = "General"

$wsAssoc.Range("A18").Value = "eventNote"
$wsAssoc.Range("A18").NumberFormat = "General"

# --- NewJobInfo: move the remembered selection to A2:A7 ---
$wsJob.Range("A2:A7").Select()

# --- emshaps90 becomes the active / selected tab, with D11 as the active cell ---
$wsAssoc.Activate()
$wsAssoc.Range("D11").Select()
